$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F6").Value = ";0;0;0;0;0;0"
$ws.Range("G6").Value = ";32;42;42;42;42;43"
$ws.Range("H6").Value = ";-750.0;-100;-100;-100;-100;-100"

$ws.Range("B6").Value = 250
$ws.Range("C6").Value = 6
$ws.Range("E6").Value = 6
